# Updates "想去人数" (F column) counts across all four worksheets
# as scraped/regenerated at commit 456a3b4 (gh-pages output refresh).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 2675   # was 2661
$ws.Range("F7").Value = 2153   # was 2096
$ws.Range("F8").Value = 1799   # was 1791
$ws.Range("F9").Value = 207   # was 206
$ws.Range("F11").Value = 2456   # was 2448
$ws.Range("F12").Value = 539   # was 537
$ws.Range("F13").Value = 227   # was 224
$ws.Range("F16").Value = 119   # was 118
$ws.Range("F17").Value = 104   # was 103
$ws.Range("F18").Value = 9088   # was 9051
$ws.Range("F19").Value = 56   # was 55
$ws.Range("F20").Value = 7048   # was 7028
$ws.Range("F21").Value = 11524   # was 11489
$ws.Range("F25").Value = 334   # was 328
$ws.Range("F26").Value = 547   # was 542
$ws.Range("F27").Value = 2540   # was 2528
$ws.Range("F28").Value = 227   # was 225
$ws.Range("F29").Value = 192   # was 191
$ws.Range("F30").Value = 2472   # was 2456
$ws.Range("F31").Value = 646   # was 635
$ws.Range("F33").Value = 4492   # was 4489
$ws.Range("F34").Value = 841   # was 819
$ws.Range("F35").Value = 342   # was 339
$ws.Range("F37").Value = 508   # was 502

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 1183   # was 1182
$ws.Range("F16").Value = 98   # was 97

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 143   # was 138

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 2675   # was 2661
$ws.Range("F9").Value = 2153   # was 2096
$ws.Range("F11").Value = 1799   # was 1791
$ws.Range("F13").Value = 207   # was 206
$ws.Range("F14").Value = 2456   # was 2448
$ws.Range("F16").Value = 539   # was 537
$ws.Range("F17").Value = 227   # was 224
$ws.Range("F20").Value = 119   # was 118
$ws.Range("F21").Value = 104   # was 103
$ws.Range("F22").Value = 9088   # was 9051
$ws.Range("F23").Value = 56   # was 55
$ws.Range("F24").Value = 7048   # was 7028
$ws.Range("F25").Value = 11524   # was 11489
$ws.Range("F29").Value = 334   # was 328
$ws.Range("F31").Value = 547   # was 542
$ws.Range("F33").Value = 2540   # was 2528
$ws.Range("F36").Value = 227   # was 225
$ws.Range("F37").Value = 192   # was 191
$ws.Range("F39").Value = 4492   # was 4489
$ws.Range("F42").Value = 98   # was 97
$ws.Range("F46").Value = 507   # was 502
